# Delete row 8 (the "MINYAK MIE SHALLOT OIL" entry under "Consume"),
# which shifts all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Delete()

# Restore the active selection to C9 (matches Excel's behavior after
# deleting a row above the previously selected cell C10).
$ws.Range("C9").Select()
